$wb = $excel.ActiveWorkbook

# --- Step 1: rename sheet tabs (swap C113<->C115, swap C111<->C112; C114 unchanged) ---
$wb.Worksheets.Item("C113").Name = "C113_tmp"
$wb.Worksheets.Item("C115").Name = "C113"
$wb.Worksheets.Item("C113_tmp").Name = "C115"

$wb.Worksheets.Item("C111").Name = "C111_tmp"
$wb.Worksheets.Item("C112").Name = "C111"
$wb.Worksheets.Item("C111_tmp").Name = "C112"

# --- Step 2: update cell contents (identified by the ORIGINAL/physical sheet name before renaming) ---

$ws = $wb.Worksheets.Item("C115")
$ws.Range("B2").Value = "Profesor:DanielL`nAsignatura:LogicaCp`nAula:4"
$ws.Range("C2").Value = "Profesor:DalianisAL`nAsignatura:AlgebraCP`nAula:3"
$ws.Range("D2").Value = "Profesor:PacoP`nAsignatura:ProgramacionCp`nAula:2"
$ws.Range("B3").Value = "Profesor:CarlaP`nAsignatura:ProgramacionCp`nAula:1"
$ws.Range("C3").Value = "Profesor:Piad`nAsignatura:Programacion`nAula:1"
$ws.Range("D3").Value = "Profesor:ErnestoA`nAsignatura:AnalisisCp`nAula:4"
$ws.Range("E3").Value = "Profesor:Idania`nAsignatura:Analisis`nAula:3"
$ws.Range("F3").Value = "Profesor:Celia`nAsignatura:Algebra`nAula:3"
$ws.Range("B4").ClearContents()
$ws.Range("D4").Value = "Profesor:Yudivian`nAsignatura:Logica`nAula:2"
$ws.Range("E4").Value = "Profesor:CristinaA`nAsignatura:AnalisisCp`nAula:1"

$ws = $wb.Worksheets.Item("C112")
$ws.Range("B2").Value = "Profesor:Idania`nAsignatura:Analisis`nAula:2"
$ws.Range("C2").Value = "Profesor:CayetanaAL`nAsignatura:AlgebraCP`nAula:5"
$ws.Range("D2").Value = "Profesor:Piad`nAsignatura:Programacion`nAula:4"
$ws.Range("E2").Value = "Profesor:CristinaA`nAsignatura:AnalisisCp`nAula:1"
$ws.Range("B3").Value = "Profesor:HectorP`nAsignatura:ProgramacionCp`nAula:3"
$ws.Range("C3").Value = "Profesor:DanielL`nAsignatura:LogicaCp`nAula:5"
$ws.Range("D3").Value = "Profesor:CarlaP`nAsignatura:ProgramacionCp`nAula:2"
$ws.Range("E3").Value = "Profesor:Celia`nAsignatura:Algebra`nAula:5"
$ws.Range("F3").Value = "Profesor:MercedesA`nAsignatura:AnalisisCp`nAula:2"
$ws.Range("C4").Value = "Profesor:Yudivian`nAsignatura:Logica`nAula:2"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "Profesor:PepeAl`nAsignatura:AlgebraCP`nAula:1"

$ws = $wb.Worksheets.Item("C114")
$ws.Range("B2").Value = "Profesor:Celia`nAsignatura:Algebra`nAula:1"
$ws.Range("C2").Value = "Profesor:Piad`nAsignatura:Programacion`nAula:4"
$ws.Range("D2").Value = "Profesor:PepeAl`nAsignatura:AlgebraCP`nAula:1"
$ws.Range("E2").Value = "Profesor:Yudivian`nAsignatura:Logica`nAula:2"
$ws.Range("B3").Value = "Profesor:CayetanaAL`nAsignatura:AlgebraCP`nAula:2"
$ws.Range("C3").Value = "Profesor:Idania`nAsignatura:Analisis`nAula:2"
$ws.Range("D3").Value = "Profesor:PacoP`nAsignatura:ProgramacionCp`nAula:5"
$ws.Range("E3").Value = "Profesor:OmarL`nAsignatura:LogicaCp`nAula:4"
$ws.Range("F3").Value = "Profesor:MercedesA`nAsignatura:AnalisisCp`nAula:2"
$ws.Range("B4").Value = "Profesor:CristinaA`nAsignatura:AnalisisCp`nAula:4"
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "Profesor:HectorP`nAsignatura:ProgramacionCp`nAula:4"

$ws = $wb.Worksheets.Item("C111")
$ws.Range("B2").Value = "Profesor:Celia`nAsignatura:Algebra`nAula:1"
$ws.Range("C2").Value = "Profesor:DanielL`nAsignatura:LogicaCp`nAula:1"
$ws.Range("D2").Value = "Profesor:ErnestoA`nAsignatura:AnalisisCp`nAula:5"
$ws.Range("E2").Value = "Profesor:CayetanaAL`nAsignatura:AlgebraCP`nAula:3"
$ws.Range("B3").Value = "Profesor:CarlaP`nAsignatura:ProgramacionCp`nAula:1"
$ws.Range("C3").Value = "Profesor:Yudivian`nAsignatura:Logica`nAula:3"
$ws.Range("D3").Value = "Profesor:MercedesA`nAsignatura:AnalisisCp`nAula:3"
$ws.Range("E3").Value = "Profesor:HectorP`nAsignatura:ProgramacionCp`nAula:1"
$ws.Range("F3").Value = "Profesor:PepeAl`nAsignatura:AlgebraCP`nAula:1"
$ws.Range("C4").Value = "Profesor:Piad`nAsignatura:Programacion`nAula:5"
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "Profesor:Idania`nAsignatura:Analisis`nAula:5"
$ws.Range("F4").ClearContents()

$ws = $wb.Worksheets.Item("C113")
$ws.Range("B2").Value = "Profesor:CarmenL`nAsignatura:LogicaCp`nAula:4"
$ws.Range("C2").Value = "Profesor:PepeAl`nAsignatura:AlgebraCP`nAula:3"
$ws.Range("D2").Value = "Profesor:ErnestoA`nAsignatura:AnalisisCp`nAula:5"
$ws.Range("E2").Value = "Profesor:Yudivian`nAsignatura:Logica`nAula:2"
$ws.Range("B3").Value = "Profesor:MercedesA`nAsignatura:AnalisisCp`nAula:5"
$ws.Range("C3").Value = "Profesor:CayetanaAL`nAsignatura:AlgebraCP`nAula:4"
$ws.Range("D3").Value = "Profesor:CarlaP`nAsignatura:ProgramacionCp`nAula:2"
$ws.Range("E3").Value = "Profesor:CarlaP`nAsignatura:ProgramacionCp`nAula:1"
$ws.Range("F3").Value = "Profesor:Celia`nAsignatura:Algebra`nAula:3"
$ws.Range("B4").Value = "Profesor:Idania`nAsignatura:Analisis`nAula:2"
$ws.Range("C4").Value = "Profesor:Piad`nAsignatura:Programacion`nAula:5"
$ws.Range("D4").ClearContents()
